$d = $word.ActiveDocument

# Remove the first paragraph in its entirety (including its paragraph
# mark), which contained the "I personally examined the patient..."
# attestation text. Range.Delete() on a paragraph's Range removes the
# text content and the paragraph mark, merging it out of the document.
$d.Paragraphs.Item(1).Range.Delete()
